$wb = $excel.ActiveWorkbook

# --- Sheet1: device_cart ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = "gugug"

# --- Sheet2: rfid_item ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A3").Value = "D792AD2D"
$ws2.Range("B3").Value = "OyVCNQgJ80lWy9HjbpvF"

$ws2.Range("A4").Value = "B8B03D1E"
$ws2.Range("B4").Value = "PXmYk7IzzsrHFMq5j70o"

$ws2.Range("A5").Value = "0A4B997F"

$ws2.Range("A6").Value = "AAAAAAAA"
$ws2.Range("B5").Value = "RMWLUuACH72OuqSPYQDk"
$ws2.Range("B6").Value = "VfgrHcX6LvHuAvkJtdgU"

$ws2.Range("A7").Value = "BBBBBBBB"
$ws2.Range("B7").Value = "YvxptylcQC7o6s7fK7H9"

$ws2.Range("A8").Value = "CCCCCCCC"
$ws2.Range("B8").Value = "oZGiQLJMymfo2Mc4KJYm"

$ws2.Range("A9").Value = "DDDDDDDD"
$ws2.Range("B9").Value = "rxRod7cigQjBK9dDmlHv"

# Resize Table2 (rfid/item) to include the two new rows
$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A1:B9"))

# Widen column A on rfid_item sheet to fit the new rfid values
$ws2.Columns.Item(1).ColumnWidth = 17

# Leave the last touched cell selected, as in the saved workbook
$null = $ws2.Range("S8").Select()

# The rfid_item sheet was the active tab when the workbook was saved
$null = $ws2.Activate()
